# Web2Assignment2Sheet.xlsx update
# - Fill in "Hours Worked" actuals for the first four tasks (Milestone 2 & 3 items)
# - Flesh out the previously-blank Milestone 5 task rows and add four new
#   milestones (6-9) with their task / details / due-date / owner / status rows
# - Re-point the selection to C26 (last thing the author clicked)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Hours worked actuals (I column) for the already-existing task rows.
#    C4 (=SUM(I8:I31)) recalculates automatically.
# ---------------------------------------------------------------------------
$ws.Range("I8").Value = 3
$ws.Range("I10").Value = 2
$ws.Range("I11").Value = 2
$ws.Range("I12").Value = 1

# ---------------------------------------------------------------------------
# 2. Milestone 5 task rows (14 & 16) - both were empty placeholder rows.
# ---------------------------------------------------------------------------

# Row 14: Login/Register/User
$ws.Range("B14").Value = "Login/Register/User"
$ws.Range("C14").Value = "Complete the function of being able to login and register. "
$ws.Range("C14").WrapText = $true
$ws.Range("D14").Value = 43796
$ws.Range("F14").Value = 43797
$ws.Range("G14").Value = "Byron, Jakub"
$ws.Range("H14").Value = "In-Progress"
$ws.Rows.Item(14).RowHeight = 51.75

# Row 16: Country/City/Photo PHP Pages + Photos in DB
$ws.Range("B16").Value = "Country/City/Photo PHP Pages + Photos in DB"
$ws.Range("B16").WrapText = $true
$ws.Range("C16").Value = "Complete the country, city, and photo php pages, as well as place photos in database."
$ws.Range("C16").WrapText = $true
$ws.Range("D16").Value = 43798
$ws.Range("F16").Value = 43799
$ws.Range("G16").Value = "Colby + Assistance"
$ws.Range("H16").Value = "In-Progress"
$ws.Rows.Item(16).RowHeight = 60

# ---------------------------------------------------------------------------
# 3. New milestone header rows (17, 19, 21, 23, 25) - copy the formatting of
#    the existing "Milestone 5" header row (15) so the grey banner style
#    (styles 16/6/22) is reproduced, then stamp in the milestone label.
# ---------------------------------------------------------------------------
$ws.Range("B15:I15").Copy()
$ws.Range("B17:I17").PasteSpecial(-4122)
$ws.Range("B19:I19").PasteSpecial(-4122)
$ws.Range("B21:I21").PasteSpecial(-4122)
$ws.Range("B23:I23").PasteSpecial(-4122)
$ws.Range("B25:I25").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("B17").Value = "Milestone 6"
$ws.Range("B19").Value = "Milestone 7"
$ws.Range("B21").Value = "Milestone 8"
$ws.Range("B23").Value = "Milestone 9"
# Row 25 stays an empty banner row (no label yet), matching the source file.

# ---------------------------------------------------------------------------
# 4. New milestone task rows (18, 20, 22, 24).
# ---------------------------------------------------------------------------

# Row 18: Search/Browse, Favourites, Uploads (Milestone 6 task)
$ws.Range("B18").Value = "Search/Browse, Favourites, Uploads"
$ws.Range("C18").Value = "Complete the search/browse funciton, favourites, as well as uploads.  "
$ws.Range("C18").WrapText = $true
$ws.Range("D18").Value = 43801
$ws.Range("F18").Value = 43802
$ws.Range("H18").Value = "Not Started"
$ws.Rows.Item(18).RowHeight = 51.75

# Row 20: Home, Profile Page, About Page (Milestone 7 task)
$ws.Range("B20").Value = "Home, Profile Page, About Page"
$ws.Range("B20").WrapText = $true
$ws.Range("C20").Value = "Home, profile, and about page should be done completely."
$ws.Range("C20").WrapText = $true
$ws.Range("D20").Value = 43803
$ws.Range("F20").Value = 43804
$ws.Range("H20").Value = "Not Started"
$ws.Rows.Item(20).RowHeight = 51.75

# Row 22: Database/Design (Milestone 8 task)
$ws.Range("B22").Value = "Database/Design"
$ws.Range("C22").Value = "Finish off the database and add finishing touches to design."
$ws.Range("C22").WrapText = $true
$ws.Range("D22").Value = 43805
$ws.Range("F22").Value = 43806
$ws.Range("H22").Value = "Not Started"
$ws.Rows.Item(22).RowHeight = 50.25

# Row 24: Testing/Submission (Milestone 9 task)
$ws.Range("B24").Value = "Testing/Submission"
$ws.Range("C24").Value = "Final Testing and Submission of assignment"
$ws.Range("C24").WrapText = $true
$ws.Range("D24").Value = 43807
$ws.Range("F24").Value = 43808
$ws.Range("H24").Value = "Not Started"
$ws.Rows.Item(24).RowHeight = 51.75

# ---------------------------------------------------------------------------
# 5. Final selection, mirroring where the author left the cursor.
# ---------------------------------------------------------------------------
$ws.Range("C26").Select()
